$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Trip Caps (New Office & Commercial Development)" row (row 9).
# This shifts the row below (EV Charging Program) up to become row 9,
# matching the target layout (A1:I9).
$ws.Rows.Item(9).Delete()

# Update the active selection to reflect the new row 9 (as in the target file).
$ws.Range("A9:XFD9").Select()
